$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date from 2021-05-17 to 2021-05-18
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.08996336095747878
$ws.Range("E2").Value = 0.01976698520748799

$ws.Range("D3").Value = 0.1053601830795341
$ws.Range("E3").Value = -0.004973291582243577

$ws.Range("D4").Value = 0.1204522999288343
$ws.Range("E4").Value = -0.008241758241758323

$ws.Range("D5").Value = 0.1419975867939534
$ws.Range("E5").Value = -0.006615941010176041

$ws.Range("D6").Value = 0.1380089300838475
$ws.Range("E6").Value = -0.004422944022114628

$ws.Range("D7").Value = 0.1484997816820724
$ws.Range("E7").Value = -0.01266654156502156

$ws.Range("D8").Value = 0.1259343409002849
$ws.Range("E8").Value = 0.003635262041805642

$ws.Range("D9").Value = 0.1297835165739946
$ws.Range("E9").Value = -0.003245288540087543

$ws.Range("E10").Value = -0.003132633939024743

# Restore sheet protection (the cells above were unlocked only to permit the edit)
$ws.Protect()
